$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 17:22"
$ws.Range("B4").Value = 7900667
$ws.Range("C4").Value = 6189
$ws.Range("D4").Value = 5065681
$ws.Range("E4").Value = 2616240
$ws.Range("G4").Value = 98
$ws.Range("H4").Value = 218746
$ws.Range("B5").Value = 6997852
$ws.Range("C5").Value = 20844
$ws.Range("D5").Value = 6003244
$ws.Range("E5").Value = 887040
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = 107568
$ws.Range("B17").Value = 479595
$ws.Range("C17").Value = 1826
$ws.Range("D17").Value = 452054
$ws.Range("E17").Value = 14269
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = 13272
$ws.Range("B20").Value = 349494
$ws.Range("C20").Value = 5724
$ws.Range("D20").Value = 238525
$ws.Range("E20").Value = 74829
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 36140
$ws.Range("B25").Value = 321392
$ws.Range("C25").Value = 914
$ws.Range("E25").Value = 38202
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9690
$ws.Range("B29").Value = 178929
$ws.Range("C29").Value = 812
$ws.Range("D29").Value = 150225
$ws.Range("E29").Value = 19112
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 9592
$ws.Range("B39").Value = 118014
$ws.Range("C39").Value = 557
$ws.Range("D39").Value = 93627
$ws.Range("E39").Value = 22220
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 2167
$ws.Range("B48").Value = 97544
$ws.Range("C48").Value = 609
$ws.Range("D48").Value = 86046
$ws.Range("E48").Value = 8133
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 3365
$ws.Range("B49").Value = 88233
$ws.Range("C49").Value = 594
$ws.Range("D49").Value = 81356
$ws.Range("E49").Value = 5253
$ws.Range("G49").Value = 8
$ws.Range("H49").Value = 1624
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 85574
$ws.Range("C51").Value = 1646
$ws.Range("D51").Value = 52803
$ws.Range("E51").Value = 30704
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 2067
$ws.Range("A52").Value = "China"
$ws.Range("B52").Value = 85536
$ws.Range("C52").Value = 15
$ws.Range("D52").Value = 80696
$ws.Range("E52").Value = 206
$ws.Range("H52").Value = 4634
$ws.Range("B58").Value = 61762
$ws.Range("C58").Value = 929
$ws.Range("D58").Value = 43982
$ws.Range("E58").Value = 16322
$ws.Range("G58").Value = 16
$ws.Range("H58").Value = 1458
$ws.Range("A73").Value = "Kenia"
$ws.Range("B73").Value = 41158
$ws.Range("C73").Value = 538
$ws.Range("D73").Value = 31876
$ws.Range("E73").Value = 8522
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 760
$ws.Range("A74").Value = "Irlanda"
$ws.Range("B74").Value = 40703
$ws.Range("D74").Value = 23364
$ws.Range("E74").Value = 15518
$ws.Range("H74").Value = 1821
$ws.Range("A85").Value = "Jordania"
$ws.Range("B85").Value = 23998
$ws.Range("C85").Value = 1235
$ws.Range("D85").Value = 6045
$ws.Range("E85").Value = 17772
$ws.Range("G85").Value = 15
$ws.Range("H85").Value = 181
$ws.Range("A86").Value = "Bulgaria"
$ws.Range("B86").Value = 23871
$ws.Range("D86").Value = 15713
$ws.Range("E86").Value = 7271
$ws.Range("H86").Value = 887
$ws.Range("A97").Value = "Albania"
$ws.Range("B97").Value = 15231
$ws.Range("C97").Value = 165
$ws.Range("D97").Value = 9406
$ws.Range("E97").Value = 5409
$ws.Range("H97").Value = 416
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 15096
$ws.Range("C98").Value = 374
$ws.Range("D98").Value = 10780
$ws.Range("E98").Value = 4161
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 155
$ws.Range("A116").Value = "Jamaica"
$ws.Range("B116").Value = 7559
$ws.Range("C116").Value = 196
$ws.Range("D116").Value = 3142
$ws.Range("E116").Value = 4279
$ws.Range("G116").Value = 6
$ws.Range("H116").Value = 138
$ws.Range("A117").Value = "Mauritania"
$ws.Range("B117").Value = 7548
$ws.Range("D117").Value = 7245
$ws.Range("E117").Value = 140
$ws.Range("H117").Value = 163
$ws.Range("B122").Value = 5948
$ws.Range("C122").Value = 5
$ws.Range("D122").Value = 5440
$ws.Range("E122").Value = 385
